$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13, shifting rows 13-21 down to 14-22.
$ws.Rows.Item(13).Insert()

# Row 10: Objetivos: text now holds the full Portuguese goal statement
# (it previously, incorrectly, held the docente text).
$ws.Range("B10").Value = "Fornecer uma visão geral sobre gestão de projetos em uma organização, suas etapas e inter-relação com as demais áreasorganizacionais e sua importância no mundo competitivo dos negócios."
$ws.Range("C10").Value = "Fornecer uma visão geral sobre gestão de projetos em uma organização, suas etapas e inter-relação com as demais áreasorganizacionais e sua importância no mundo competitivo dos negócios."

# Row 13 (new row): only B/C are used (Docentes responsáveis value,
# previously stuck on row 10); column A stays empty for this row.
$ws.Range("A13").Clear()
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840917 - Fabricio Maciel Gomes"
$ws.Range("C13").Value = "5840917 - Fabricio Maciel Gomes"

# Row 14: Programa resumido: correct Portuguese summary (was "Semestral").
$ws.Range("B14").Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."
$ws.Range("C14").Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."

# Row 16: Programa: same Portuguese summary (was a stray date).
$ws.Range("B16").Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."
$ws.Range("C16").Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."

# Row 19: Método: correct method text (was the docente text).
$ws.Range("B19").Value = "Aulas expositivas. Trabalhos em grupo. Seminários. Palestras. Exercícios em sala de aula."
$ws.Range("C19").Value = "Aulas expositivas. Trabalhos em grupo. Seminários. Palestras. Exercícios em sala de aula."

# Row 20: Critério: correct grading criteria text (was the método text).
$ws.Range("B20").Value = "Duas Provas com peso de 30% cada uma. Trabalhos em sala de aula com peso de 20% e Trabalho final com peso de 20%"
$ws.Range("C20").Value = "Duas Provas com peso de 30% cada uma. Trabalhos em sala de aula com peso de 20% e Trabalho final com peso de 20%"

# Row 21: Norma de recuperação: "Prova única" (was the critério text).
$ws.Range("B21").Value = "Prova única"
$ws.Range("C21").Value = "Prova única"

# Row 22: Bibliografia: full bibliography text (was "Prova única").
$ws.Range("B22").Value = "1. PMBOK. Um Guia Do Conhecimento em Gerenciamento de projetos. 5 ed. Project Management Institute. 20122. CLAUSING, D. Total quality development a step by step guide to world class concurrent engineering. New York: ASME Press,1994.3. MEREDITH, J R; MANTEL, S J; WILEY, J. Project Management: a managerial approach. 1995.4. MAXIMIANO, A . C. Administração de projetos, Atlas: São Paulo, 1997.5. SHTUB, A BARD J. F. e GLOBERSON S. Project management, Prentice hall, 1994."
$ws.Range("C22").Value = "1. PMBOK. Um Guia Do Conhecimento em Gerenciamento de projetos. 5 ed. Project Management Institute. 20122. CLAUSING, D. Total quality development a step by step guide to world class concurrent engineering. New York: ASME Press,1994.3. MEREDITH, J R; MANTEL, S J; WILEY, J. Project Management: a managerial approach. 1995.4. MAXIMIANO, A . C. Administração de projetos, Atlas: São Paulo, 1997.5. SHTUB, A BARD J. F. e GLOBERSON S. Project management, Prentice hall, 1994."

Write-Output "Edit applied successfully"
